$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Experimental: true -> (cleared)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# (use formula+copy/paste-values trick so Excel doesn't auto-convert the
#  text into a date serial number, keeping it a plain shared string)
$dateCell = $ws.Range("B8")
$dateCell.Formula = '="2025-11-18"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)
